# Apply scheduled runner updates to market-price-derived columns (H-N)
# across multiple sheets/rows, per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 2227.55
$ws.Range("I33").Value = 210.33333
$ws.Range("J33").Value = 8279.200000000001
$ws.Range("K33").Value = 210.33333
$ws.Range("L33").Value = 8279.200000000001
$ws.Range("M33").Value = 18.66667000000001
$ws.Range("N33").Value = -8737.200000000001
# Row 43
$ws.Range("H43").Value = 10859.375
$ws.Range("I43").Value = 2456.3333
$ws.Range("J43").Value = 15901.2
$ws.Range("K43").Value = 2456.3333
$ws.Range("L43").Value = 15901.2
$ws.Range("M43").Value = -2387.3333
$ws.Range("N43").Value = -16039.2
# Row 75
$ws.Range("H75").Value = 30285
$ws.Range("I75").Value = 30285
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 30285
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -29349
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 30285
$ws.Range("I78").Value = 30285
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 90855
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -86175
$ws.Range("N78").ClearContents()
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2208.745
$ws.Range("I32").Value = 1877.7021
$ws.Range("K32").Value = 1877.7021
$ws.Range("M32").Value = -1590.7021
# Row 74
$ws.Range("H74").Value = 9261838
$ws.Range("I74").Value = 9806387
$ws.Range("K74").Value = 9806387
$ws.Range("M74").Value = -9805513
# Row 77
$ws.Range("H77").Value = 9261838
$ws.Range("I77").Value = 9806387
$ws.Range("K77").Value = 49031935
$ws.Range("M77").Value = -49027567

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 5978.5
$ws.Range("I64").Value = 2032.6666
$ws.Range("J64").Value = 9924.333000000001
$ws.Range("K64").Value = 2032.6666
$ws.Range("L64").Value = 9924.333000000001
$ws.Range("M64").Value = -1807.6666
$ws.Range("N64").Value = -10374.333
# Row 67
$ws.Range("H67").Value = 5978.5
$ws.Range("I67").Value = 2032.6666
$ws.Range("J67").Value = 9924.333000000001
$ws.Range("K67").Value = 2032.6666
$ws.Range("L67").Value = 9924.333000000001
$ws.Range("M67").Value = -1252.6666
$ws.Range("N67").Value = -11484.333
# Row 94
$ws.Range("H94").Value = 1240.5454
$ws.Range("I94").Value = 1240.5454
$ws.Range("K94").Value = 1240.5454
$ws.Range("M94").Value = -789.5454
# Row 138
$ws.Range("H138").Value = 65023.582
$ws.Range("J138").Value = 65023.582
$ws.Range("L138").Value = 65023.582
$ws.Range("N138").Value = -75303.58199999999

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 8338.546
$ws.Range("I58").Value = 4598
$ws.Range("J58").Value = 11455.667
$ws.Range("K58").Value = 4598
$ws.Range("L58").Value = 11455.667
$ws.Range("M58").Value = -4395
$ws.Range("N58").Value = -11861.667
# Row 62
$ws.Range("H62").Value = 10899.5
$ws.Range("I62").Value = 8373.625
$ws.Range("J62").Value = 21003
$ws.Range("K62").Value = 8373.625
$ws.Range("L62").Value = 21003
$ws.Range("M62").Value = -7749.625
$ws.Range("N62").Value = -22251
# Row 65
$ws.Range("H65").Value = 10899.5
$ws.Range("I65").Value = 8373.625
$ws.Range("J65").Value = 21003
$ws.Range("K65").Value = 41868.125
$ws.Range("L65").Value = 105015
$ws.Range("M65").Value = -38748.125
$ws.Range("N65").Value = -111255
# Row 134
$ws.Range("H134").Value = 2366.8647
$ws.Range("I134").Value = 1780.6061
$ws.Range("K134").Value = 5341.8183
$ws.Range("M134").Value = -2806.8183
# Row 136
$ws.Range("H136").Value = 8338.546
$ws.Range("I136").Value = 4598
$ws.Range("J136").Value = 11455.667
$ws.Range("K136").Value = 13794
$ws.Range("L136").Value = 34367.001
$ws.Range("M136").Value = -11244
$ws.Range("N136").Value = -39467.001
# Row 141
$ws.Range("H141").Value = 233996
$ws.Range("J141").Value = 233996
$ws.Range("L141").Value = 233996
$ws.Range("N141").Value = -244356

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 7437.591
$ws.Range("I5").Value = 942.6667
$ws.Range("K5").Value = 2828.0001
$ws.Range("M5").Value = -2716.0001
# Row 132
$ws.Range("H132").Value = 3399.5518
$ws.Range("J132").Value = 3213.158
$ws.Range("L132").Value = 28918.422
$ws.Range("N132").Value = -33978.422
# Row 135
$ws.Range("H135").Value = 7437.591
$ws.Range("I135").Value = 942.6667
$ws.Range("K135").Value = 8484.0003
$ws.Range("M135").Value = -5949.0003
# Row 136
$ws.Range("H136").Value = 1588.5
$ws.Range("I136").Value = 1227.909
$ws.Range("K136").Value = 3683.727
$ws.Range("M136").Value = 1416.273
# Row 137
$ws.Range("H137").Value = 2216.5
$ws.Range("I137").Value = 1400
$ws.Range("J137").Value = 3033
$ws.Range("K137").Value = 4200
$ws.Range("L137").Value = 9099
$ws.Range("M137").Value = 900
$ws.Range("N137").Value = -19299
# Row 139
$ws.Range("H139").Value = 3399.476
$ws.Range("I139").Value = 1287.1875
$ws.Range("K139").Value = 3861.5625
$ws.Range("M139").Value = 1278.4375

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
# Row 126
$ws.Range("H126").Value = 5706.7295
$ws.Range("I126").Value = 2632.842
$ws.Range("J126").Value = 8951.388999999999
$ws.Range("K126").Value = 7898.526
$ws.Range("L126").Value = 26854.167
$ws.Range("M126").Value = -5428.526
$ws.Range("N126").Value = -31794.167
# Row 132
$ws.Range("H132").Value = 7157.393
$ws.Range("I132").Value = 1814.5
$ws.Range("J132").Value = 10125.667
$ws.Range("K132").Value = 5443.5
$ws.Range("L132").Value = 30377.001
$ws.Range("M132").Value = -2913.5
$ws.Range("N132").Value = -35437.001

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 2941629.8
$ws.Range("I55").Value = 4545900.5
$ws.Range("K55").Value = 4545900.5
$ws.Range("M55").Value = -4545727.5
# Row 61
$ws.Range("H61").Value = 3313.25
$ws.Range("I61").Value = 2334.8333
$ws.Range("J61").Value = 6248.5
$ws.Range("K61").Value = 2334.8333
$ws.Range("L61").Value = 6248.5
$ws.Range("M61").Value = -2132.8333
$ws.Range("N61").Value = -6652.5
# Row 113
$ws.Range("H113").Value = 3313.25
$ws.Range("I113").Value = 2334.8333
$ws.Range("J113").Value = 6248.5
$ws.Range("K113").Value = 2334.8333
$ws.Range("L113").Value = 6248.5
$ws.Range("M113").Value = -164.8332999999998
$ws.Range("N113").Value = -10588.5
# Row 136
$ws.Range("H136").Value = 7614.6743
$ws.Range("I136").Value = 5389.0625
$ws.Range("J136").Value = 14089.182
$ws.Range("K136").Value = 16167.1875
$ws.Range("L136").Value = 42267.546
$ws.Range("M136").Value = -13617.1875
$ws.Range("N136").Value = -47367.546

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2766
$ws.Range("I96").Value = 3250
$ws.Range("J96").Value = 1798
$ws.Range("K96").Value = 3250
$ws.Range("L96").Value = 1798
$ws.Range("M96").Value = -1877
$ws.Range("N96").Value = -4544
# Row 135
$ws.Range("H135").Value = 69249.25
$ws.Range("J135").Value = 69249.25
$ws.Range("L135").Value = 69249.25
$ws.Range("N135").Value = -79389.25
# Row 136
$ws.Range("H136").Value = 3641.3
$ws.Range("I136").Value = 2109.077
$ws.Range("K136").Value = 6327.231000000001
$ws.Range("M136").Value = -3777.231000000001
